$wb = $excel.ActiveWorkbook

# --- Sheet 1: Training Dashboard -------------------------------------------
$ws1 = $wb.Worksheets.Item("Training Dashboard")

# Row 3 (LOTO SOPs) - refreshed "period to expire" + "last update" date
$ws1.Range("H3").Value = 85
$ws1.Range("I3").NumberFormat = "@"
$ws1.Range("I3").Value = "16-Sep-2025"

# Row 4 (Endangered by Electricity) - refreshed "period to expire" + "last update" date
$ws1.Range("H4").Value = -50
$ws1.Range("I4").NumberFormat = "@"
$ws1.Range("I4").Value = "16-Sep-2025"

# --- Sheet 2: Exam Dashboard -------------------------------------------------
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

# --- Header styling -----------------------------------------------------------
# The dashboard title and the column-header row now share a single bold,
# white font (rendered on top of the dark-blue header fill) instead of two
# separate bold fonts, so bring the title down to the header's font size
# and recolor both to white.
$ws1.Range("A2:K2").Font.Color = 16777215
$ws1.Range("A1").Font.Color = 16777215
$ws1.Range("A1").Font.Size = 11

$ws2.Range("A2:G2").Font.Color = 16777215
$ws2.Range("A1").Font.Color = 16777215
$ws2.Range("A1").Font.Size = 11
